# Add a "수집일" (collection date) column to the scraped-listing sheet.
# Mirrors the source workbook's naming convention (df_list_20241128.xlsx):
# every row gets stamped with the collection date 2024-11-28 (serial 45624).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last populated data row in column A (header is row 1).
$lastRow = $ws.Cells.Item($ws.Rows.Count(), 1).End(-4162).Row()

# --- Header cell F1 -------------------------------------------------
$headerCell = $ws.Cells.Item(1, 6)
$headerCell.Value = "수집일"
$headerCell.Font.Bold = $true
$headerCell.HorizontalAlignment = -4108
$headerCell.VerticalAlignment = -4160
$headerCell.Borders.Item(7).LineStyle = 1
$headerCell.Borders.Item(10).LineStyle = 1

# --- Data cells F2:F<lastRow> ---------------------------------------
$dataRange = $ws.Range($ws.Cells.Item(2, 6), $ws.Cells.Item($lastRow, 6))
$dataRange.NumberFormat = "m/d/yyyy"
$dataRange.Value = 45624
